$d = $word.ActiveDocument
$divide = [char]0x00F7

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2026-02-17 Tuesday" "2026-02-18 Wednesday"

Replace-Text "826${divide}6=" "750${divide}5="
Replace-Text "933${divide}9=" "847${divide}2="
Replace-Text "128${divide}3=" "626${divide}5="
Replace-Text "264${divide}4=" "340${divide}7="
Replace-Text "738${divide}8=" "915${divide}2="
Replace-Text "865${divide}8=" "770${divide}4="
Replace-Text "959${divide}9=" "941${divide}6="
Replace-Text "743${divide}6=" "105${divide}7="
Replace-Text "974${divide}8=" "337${divide}9="
Replace-Text "224${divide}3=" "862${divide}6="
Replace-Text "540${divide}3=" "443${divide}8="
Replace-Text "929${divide}5=" "207${divide}7="
Replace-Text "955${divide}9=" "472${divide}5="
Replace-Text "133${divide}4=" "372${divide}6="
Replace-Text "797${divide}7=" "484${divide}3="
Replace-Text "736${divide}6=" "534${divide}9="
Replace-Text "331${divide}7=" "687${divide}9="
Replace-Text "383${divide}9=" "877${divide}3="
Replace-Text "973${divide}5=" "423${divide}2="
Replace-Text "286${divide}9=" "874${divide}8="
Replace-Text "281${divide}6=" "568${divide}7="
Replace-Text "240${divide}7=" "244${divide}4="
Replace-Text "511${divide}3=" "633${divide}7="
Replace-Text "836${divide}2=" "580${divide}9="
Replace-Text "530${divide}5=" "544${divide}2="
